$d = $word.ActiveDocument

# The last paragraph of the poem ends in "... ждёт…". Two new closing
# lines need to be appended as their own new paragraphs, right before the
# closing section properties:
#   "А мне старушка теплые варенки плетет,"
#   "Но не учла она что я побольше буду"
#
# Using two separate Find/Replace passes (rather than one combined
# replacement, and rather than Range.InsertParagraphAfter) keeps each new
# paragraph's run free of the inherited/"remembered" character formatting
# (the white font color used throughout the rest of the poem), matching
# the plain, unformatted runs introduced by the edit.

$d.Content.Find.Execute(
    " ждёт…", $true, $false, $false, $false, $false, $true, 1, $false,
    " ждёт…^pА мне старушка теплые варенки плетет,", 2) | Out-Null

$d.Content.Find.Execute(
    "А мне старушка теплые варенки плетет,", $true, $false, $false, $false, $false, $true, 1, $false,
    "А мне старушка теплые варенки плетет,^pНо не учла она что я побольше буду", 2) | Out-Null
